# إضافة حدث جديد في Card20
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card20")

# Row 19: fill in the previously-blank "nan" placeholder cells (B..K, M)
$ws.Range("B19").Value = "nan"
$ws.Range("C19").Value = "nan"
$ws.Range("D19").Value = "nan"
$ws.Range("E19").Value = "nan"
$ws.Range("F19").Value = "nan"
$ws.Range("G19").Value = "nan"
$ws.Range("H19").Value = "nan"
$ws.Range("I19").Value = "nan"
$ws.Range("J19").Value = "nan"
$ws.Range("K19").Value = "nan"
$ws.Range("M19").Value = "nan"

# Row 20: new service event entry
# Copy A19 ("20", stored as text) down into A20 so it stays text, not a number
$ws.Range("A19").Copy($ws.Range("A20"))
$ws.Range("L20").Value = "3\2\2024"
$ws.Range("N20").Value = "تم تغيير سلك الدوفر"
$ws.Range("O20").Value = "الخبير"
